# Pedidos.xlsx update — "Add files via upload" / "atualizar Pedidos"
# Fills in the previously-blank rows 29-49 with newly received order data
# and pushes the previous row-29 entry down to row 50 (matching the rest
# of the already-populated list below it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A = Remessa, Column B = Material, Column C = Quantidade
$data = @(
    @{ Row = 29; A = "80266517"; B = "10377-ARI-I"; C = 1 },
    @{ Row = 30; A = "80266519"; B = "15386-DLO-I"; C = 1 },
    @{ Row = 31; A = "80266522"; B = "10361-ARI-I"; C = 1 },
    @{ Row = 32; A = "80266523"; B = "10000-LDG-I"; C = 1 },
    @{ Row = 33; A = "80266524"; B = "10533-DLO-I"; C = 3 },
    @{ Row = 34; A = "80266525"; B = "10020-ARI-I"; C = 1 },
    @{ Row = 35; A = "80266525"; B = "10388-ARI-I"; C = 1 },
    @{ Row = 36; A = "80266526"; B = "10383-ARI-I"; C = 1 },
    @{ Row = 37; A = "80266527"; B = "10000-LDG-I"; C = 1 },
    @{ Row = 38; A = "80266529"; B = "33664-ATE-I"; C = 1 },
    @{ Row = 39; A = "80266530"; B = "10125-ARI-I"; C = 1 },
    @{ Row = 40; A = "80266534"; B = "33642-ATE-I"; C = 3 },
    @{ Row = 41; A = "80266534"; B = "33652-ATE-I"; C = 2 },
    @{ Row = 42; A = "80266534"; B = "33382-ATE-I"; C = 8 },
    @{ Row = 43; A = "80266534"; B = "33380-ATE-I"; C = 40 },
    @{ Row = 44; A = "80266534"; B = "23495-GPB-I"; C = 10 },
    @{ Row = 45; A = "80266534"; B = "23432-GPB-I"; C = 10 },
    @{ Row = 46; A = "80266534"; B = "10000-MBY-I"; C = 1 },
    @{ Row = 47; A = "80266534"; B = "21501-NZX-I"; C = 1 },
    @{ Row = 48; A = "80266537"; B = "10645-ARI-I"; C = 2 },
    @{ Row = 49; A = "80266537"; B = "10486-ARI-I"; C = 2 },
    @{ Row = 50; A = "84004823"; B = "15211-DLO-I"; C = 1 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
}

# Restore the cursor/selection left behind when the upload was made
[void]$ws.Range("B12").Select()
